$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column I: hw8 ---

# Header (I1): copy style from H1 (bold header style) then set the text.
$ws.Range("H1").Copy() | Out-Null
$ws.Range("I1").PasteSpecial(-4122) | Out-Null
$ws.Range("I1").Value = "hw8"

# hw8 has only been partially graded so far - only some students have a
# score recorded yet (rows 3,4,5,7,8,9,10,11,12). Copy the plain data
# style from H3 into just those specific cells (not the whole column)
# so rows 2, 6, 13, 14, 15, 16 are left untouched.
$targetRows = @(3, 4, 5, 7, 8, 9, 10, 11, 12)
foreach ($r in $targetRows) {
    $ws.Range("H$r").Copy() | Out-Null
    $ws.Range("I$r").PasteSpecial(-4122) | Out-Null
}

$ws.Range("I3").Value = 97
$ws.Range("I4").Value = 99
$ws.Range("I5").Value = 98
$ws.Range("I7").Value = 97
$ws.Range("I8").Value = 100
$ws.Range("I9").Value = 99
$ws.Range("I10").Value = 100
$ws.Range("I11").Value = 99
$ws.Range("I12").Value = 99

# Match column width to the other grade columns (C:H).
$ws.Range("I1").ColumnWidth = $ws.Range("C1").ColumnWidth

# Final selection lands on I13, matching the author's last-clicked cell.
$ws.Range("I13").Select() | Out-Null
